$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data revisions coming from the source ("MV -datos-") ---

# Row 73 (period with index 78 in Serie list)
$ws.Range("B73").Value = 8417.129999999999
$ws.Range("E73").Value = 6084.17

# Row 129
$ws.Range("C129").Value = 252.53
$ws.Range("D129").Value = 1468.63

# Row 130
$ws.Range("C130").Value = 248
$ws.Range("D130").Value = 1568.94

# Row 131
$ws.Range("C131").Value = 233.68
$ws.Range("D131").Value = 1589.03

# Row 132
$ws.Range("C132").Value = 236.66
$ws.Range("D132").Value = 1611.22

# --- Append a new monthly row (01-09-2021) at the bottom of the table ---

# Build the period label as genuine text (not an auto-converted date serial)
# by using a TEXT() formula in a scratch cell, then paste only the resulting
# value into the target cell. This mirrors how the source text labels in
# column A (e.g. "01-08-2021") are stored as plain shared strings.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=TEXT(""01-09-2021"",""@"")"
$scratch.Copy()
$ws.Range("A140").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B140").Value = 8345.24
$ws.Range("C140").Value = 254.08
$ws.Range("D140").Value = 1791.54
$ws.Range("E140").Value = 6028.97
$ws.Range("F140").Value = 212.04
$ws.Range("G140").Value = 58.61
